$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C10").Value = 8621
$ws.Range("C11:C13").Value = 8492
$ws.Range("C14:C19").Value = 8035
$ws.Range("C20:C59").Value = 7765
$ws.Range("C60:C65").Value = 7318
$ws.Range("C66:C83").Value = 7310
